# Add a seventh testcase: a new worksheet "Лист6" at the end of the
# workbook, built as a copy of "Лист5" (same 4-row merged layout), with
# its own test-case content (Function isInBound in maze).

$wb = $excel.ActiveWorkbook

# --- 1. Leave the previously-active sheet (Лист5) with the selection it
#        had when the user tabbed away to work on the new sheet.
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("G2:G5").Select()

# --- 2. Duplicate Лист5 to the end of the workbook and rename it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$ws6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6.Name = "Лист6"

# --- 3. Update the ID number for the new test case.
$ws6.Range("A2").Value = 6

# --- 4. Test Case name (B2, merged B2:B5).
$ws6.Range("B2").Value = "Function `nisInBound in maze"

# --- 5. Testing Steps (split across the C2:C3 / C4:C5 merges).
$ws6.Range("C2").Value = "1.Check if the function takes the inputed width and height"
$ws6.Range("C4").Value = '2.Check the weather the "if" condition statements in the function work'

# --- 6. Expected results — split D2:D5 into D2:D3 / D4:D5 (it used to be
#        one single D2:D5 merge on Лист5).
$ws6.Range("D2:D5").UnMerge()
$ws6.Range("D2:D3").Merge()
$ws6.Range("D4:D5").Merge()
$ws6.Range("D2").Value = "1.The function takes the with and height"
$ws6.Range("D4").Value = '2. The "if" condition statements work and a result is received'

# --- 7. Status (E2:E5) stays "Passed" — untouched, already copied.

# --- 8. Actual Result — split F2:F5 into F2:F3 / F4:F5.
$ws6.Range("F2:F5").UnMerge()
$ws6.Range("F2:F3").Merge()
$ws6.Range("F4:F5").Merge()
$ws6.Range("F2").Value = "1.The function successfully takes the with and height"
$ws6.Range("F4").Value = '2.The "if" condition statements works and a the function returns "true" or "false"'

# --- 9. Comments (G2:G5) stays as-is — untouched, already copied.

# --- 9b. Лист5 carried a couple of stray formatted-but-empty cells on row
#         6 (B6/C6) that don't belong on the fresh test-case sheet.
$ws6.Range("B6:C6").Clear()

# --- 10. Row heights to fit the new wrapped text.
$ws6.Rows.Item(2).RowHeight = 51
$ws6.Rows.Item(3).AutoFit()
$ws6.Rows.Item(4).RowHeight = 55.5
$ws6.Rows.Item(5).RowHeight = 15.75

# --- 11. Column widths, re-fit for the new content.
$ws6.Columns.Item(2).ColumnWidth = 19.5703125
$ws6.Columns.Item(3).ColumnWidth = 19.7109375
$ws6.Columns.Item(4).ColumnWidth = 25.140625
$ws6.Columns.Item(5).ColumnWidth = 10.42578125
$ws6.Columns.Item(6).ColumnWidth = 23.28515625
$ws6.Columns.Item(7).ColumnWidth = 19.140625

# --- 12. Selection / active cell on the new sheet, which also makes it
#         the active tab of the workbook.
$ws6.Range("D9").Select()
